$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a purely-numeric-looking value into a cell while keeping it
# stored as TEXT (matching how this workbook's existing numeric-looking
# columns - Age, Phone, Years - are stored as shared strings, not numbers).
# A plain `.Value = "8"` assignment gets auto-converted to a real number by
# Excel, so instead we stage the literal text in a scratch cell via a
# quoted-string formula, copy it, and paste-special (values only) into the
# destination - this carries the Text type across without leaving any
# formula or extra number-format/style behind. The scratch cell is cleared
# immediately after each use.
function Set-TextValue($range, [string]$text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = "=""" + $text + """"
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.Value = ""
}

# Row 1 (existing candidate Marco Simoncelli): only the last column changes,
# from 10 to 8
Set-TextValue $ws.Range("G1") "8"

# Row 2: previous candidate (Simona Laurenzi) is replaced with the new one
# extracted from the incoming PDF attachment (Iacopo De Palatis)
$ws.Range("A2").Value = "Iacopo"
$ws.Range("B2").Value = "De Palatis"
Set-TextValue $ws.Range("C2") "25"
$ws.Range("D2").Value = "iacopo@gmail.com"
Set-TextValue $ws.Range("E2") "3232854389"
$ws.Range("F2").Value = "Full-Stack"
Set-TextValue $ws.Range("G2") "10"
